$d = $word.ActiveDocument

# Common XML package wrapper used for Range.InsertXML calls below.
# InsertXML REPLACES the *entire* contents of the block (paragraph) that
# the target Range lives in, so each fragment below is a complete,
# self-contained replacement <w:p> carrying all of the original
# paragraph-level attributes (paraId/textId/rsid.../pPr) plus the new run
# content, to avoid losing anything that isn't supposed to change.
$pkgOpen = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$tbl = $d.Tables.Item(1)

# --- 1) Row 2 (Vangala Ruchitha) / Skills cell: split "HTML, Javascript, "
#        into three runs, wrapping "Javascript" in spellcheck proofErr tags.
$skillsCell1 = $tbl.Cell(2, 2)
$p1 = '<w:p w14:paraId="4597E854" w14:textId="04FE474F" w:rsidR="005C2426" w:rsidRPr="00BE429F" w:rsidRDefault="00BE429F">' `
  + '<w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr>' `
  + '<w:r w:rsidRPr="00BE429F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">HTML, </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r w:rsidRPr="00BE429F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Javascript</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r w:rsidRPr="00BE429F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r w:rsidRPr="00BE429F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Angular</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r w:rsidRPr="00BE429F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> JS, Java, MySQ</w:t></w:r>' `
  + '<w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>L</w:t></w:r>' `
  + '</w:p>'
$skillsCell1.Range.InsertXML($pkgOpen + $p1 + $pkgClose)

# --- 2) Row 3 (first blank row) / Name cell: add "Nandhini Kasukurthi",
#        each word wrapped in spellcheck proofErr tags.
$nameCell2 = $tbl.Cell(3, 1)
$p2 = '<w:p w14:paraId="32F29070" w14:textId="77777777" w:rsidR="005C2426" w:rsidRPr="00BE429F" w:rsidRDefault="005C2426">' `
  + '<w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Nandhini</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Kasukurthi</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '</w:p>'
$nameCell2.Range.InsertXML($pkgOpen + $p2 + $pkgClose)

# --- 3) Row 3 / Skills cell: add "Design and implementation" plus the
#        _GoBack bookmark that Word drops at the last edit position.
$skillsCell2 = $tbl.Cell(3, 2)
$p3 = '<w:p w14:paraId="1E9BEB56" w14:textId="77777777" w:rsidR="005C2426" w:rsidRPr="00BE429F" w:rsidRDefault="005C2426">' `
  + '<w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Design and implementation</w:t></w:r>' `
  + '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' `
  + '<w:bookmarkEnd w:id="0"/>' `
  + '</w:p>'
$skillsCell2.Range.InsertXML($pkgOpen + $p3 + $pkgClose)
